$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text would otherwise be auto-coerced
# into a Number by Excel (stripping leading/trailing zeros, changing precision).
# Force them to Text format first so the literal string is preserved exactly,
# matching the original workbook which stores these as inline strings.
$textCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D12","D14","D15","D16","D18","D19","D20","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D45","D46","D47","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated crypto market values scraped by the GitHub Actions job.
$ws.Range("D2").Value = "30.313.27"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.929.37"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "0.7472"
$ws.Range("E5").Value = "  +4.13%  "
$ws.Range("D6").Value = "243.80"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.3171"
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("D9").Value = "27.47"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").Value = "0.07112"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").Value = "0.7807"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").Value = "0.08049"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "1.921.63"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").Value = "5.400"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "93.29"
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("D16").Value = "14.57"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "30.335.07"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "6.034"
$ws.Range("E18").Value = "  +4.83%  "
$ws.Range("D19").Value = "251.83"
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("D20").Value = "0.000007919"
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.231.58"
$ws.Range("E21").Value = "  +2.25%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "6.666"
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("D25").Value = "9.583"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").Value = "165.48"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").Value = "19.07"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "0.1290"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").Value = "2.181"
$ws.Range("E29").Value = "  -3.65%  "
$ws.Range("D30").Value = "1.564"
$ws.Range("E30").Value = "  +2.23%  "
$ws.Range("D31").Value = "1.360"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").Value = "4.426"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "4.141"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "0.05248"
$ws.Range("E34").Value = "  +2.38%  "
$ws.Range("D35").Value = "1.316"
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("D36").Value = "0.7579"
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("D37").Value = "2.769"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "0.01952"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "6.519"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D41").Value = "77.53"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").Value = "0.4521"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "1.970"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("D44").Value = "0.8432"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "10.03"
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("D47").Value = "7.697"
$ws.Range("E47").Value = "  +3.20%  "
$ws.Range("D48").Value = "101.63"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").Value = "2.111.43"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").Value = "37.88"
$ws.Range("E50").Value = "  +2.78%  "
$ws.Range("D51").Value = "0.1217"
$ws.Range("E51").Value = "  +6.71%  "
